$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.808.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.384.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.368.76'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.173'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.16%  '
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.56'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000278'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.925.29'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.379.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '65.749.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.85'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '467.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.26%  '
$ws.Range("E23").Value = '  -3.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.60'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.19%  '
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.25%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.90%  '
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.03%  '
$ws.Range("E29").Value = '  -1.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.84%  '
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '577.13'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '61.96'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.06%  '
$ws.Range("E35").Value = '  -0.74%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.142'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.89'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.376'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0739'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.095.93'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.83'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0417'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.44'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.34%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.997'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '139.87'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.01%  '
$ws.Range("E50").Value = '  -1.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.96%  '
